# Update the tab-stop layout for the four "placeholder" signature blocks
# (the ones containing the "date_short" placeholder and sharing the
# 3402/3686/7655/8222 twip tab pattern). Adds a new left tab at 3119 and
# shifts the other three stops to 3402/7371/8080.
$d = $word.ActiveDocument

$oldTabsTwips = @(3402, 3686, 7655, 8222)
$newTabsTwips = @(3119, 3402, 7371, 8080)

$search = $d.Content
$search.Start = 0
$search.End = $d.Content.End

$guard = 0
while ($search.Find.Execute("date_short", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $guard = $guard + 1
    if ($guard -gt 100) { break }

    $para = $search.Paragraphs(1)
    $tabStops = $para.Range.ParagraphFormat.TabStops

    $count = $tabStops.Count
    $positions = @()
    for ($j = 1; $j -le [Math]::Min($count, $oldTabsTwips.Length); $j++) {
        $positions += [Math]::Round($tabStops.Item($j).Position * 20)
    }

    $isMatch = $true
    if ($positions.Length -ne $oldTabsTwips.Length) {
        $isMatch = $false
    } else {
        for ($k = 0; $k -lt $oldTabsTwips.Length; $k++) {
            if ($positions[$k] -ne $oldTabsTwips[$k]) {
                $isMatch = $false
            }
        }
    }

    if ($isMatch) {
        $tabStops.ClearAll()
        foreach ($posTwips in $newTabsTwips) {
            $tabStops.Add($posTwips / 20)
        }
    }

    # Move past this hit so the next Find.Execute looks further down the doc.
    $search.Collapse(0)
    $search.End = $d.Content.End
}
